$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear old row 2 formulas/data entirely ---
$ws.Range("A2:J2").Clear()

# --- Clear old headers beyond column H (I1, J1 no longer used) ---
$ws.Range("I1:J1").Clear()

# --- Update header row text (columns shift meaning/position) ---
$ws.Range("A1").Value = "Stock ID"
$ws.Range("B1").Value = "Price"
$ws.Range("C1").Value = "Yield"
$ws.Range("D1").Value = "Annual Yield"
$ws.Range("E1").Value = "`$price/`$annual"
$ws.Range("F1").Value = "Annual yield for `$1k"
$ws.Range("G1").Value = "Updated:"
$ws.Range("H1").Value = "2019-07-30 16:35:34.485624"

# --- Apply currency number format to columns B, D, E, F ---
$ws.Columns.Item(2).NumberFormat = "`"`$`"#,##0.00"
$ws.Columns.Item(4).NumberFormat = "`"`$`"#,##0.00"
$ws.Columns.Item(5).NumberFormat = "`"`$`"#,##0.00"
$ws.Columns.Item(6).NumberFormat = "`"`$`"#,##0.00"

# --- Column widths approximating autofit result on the new headers ---
$ws.Columns.Item(1).ColumnWidth = 6.3333333333333333
$ws.Columns.Item(2).ColumnWidth = 3.6666666666666665
$ws.Columns.Item(3).ColumnWidth = 3.6666666666666665
$ws.Columns.Item(4).ColumnWidth = 9.5
$ws.Columns.Item(5).ColumnWidth = 11.833333333333334
$ws.Columns.Item(6).ColumnWidth = 15.5
$ws.Columns.Item(7).ColumnWidth = 7.1666666666666667
$ws.Columns.Item(8).ColumnWidth = 23.166666666666668

# --- Sheet view: select G5, mark tab selected ---
$ws.Range("G5").Select()
$ws.Activate()

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1
